$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Paragraph 1: "Link of the code: <hyperlink>" + trailing line break
# -------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1FullRange = $d.Range($p1.Range.Start, $p1.Range.End)
$p1FullRange.Delete()

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertParagraphBefore()

$p1 = $d.Paragraphs(1)
$newUrl1 = "https://drive.google.com/file/d/1_7X-JiXTV_ttVsgLB5lncuglpa1XcR2E/view?usp=sharing"
$p1Content = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$p1Content.Text = "Link of the code: " + $newUrl1

$p1 = $d.Paragraphs(1)
$urlRange = $d.Range($p1.Range.Start, $p1.Range.End)
$urlRange.Find.Execute($newUrl1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($urlRange, $newUrl1, "", "", $newUrl1) | Out-Null

$p1 = $d.Paragraphs(1)
$lineBreakRange = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$lineBreakRange.InsertBreak(6)

# -------------------------------------------------------------------
# Paragraph 3: "Link of the Folder: <new folder url>"
# -------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$newUrl2 = "https://drive.google.com/drive/folders/1fpsYh5pCdFFl-Z5xXKGUmyRX1Yhc3Z-c?usp=sharing"
$p3.Range.Find.Execute("https://drive.google.com/drive/folders/1vwaEeWJkjxfodMKgv-Aj1n5NPkgzijz4", $false, $false, $false, $false, $false, $true, 1, $false, $newUrl2, 2) | Out-Null

# -------------------------------------------------------------------
# Bulleted list formatting for paragraphs 1 and 3 (sharing the same
# numbering list), paragraph 2 (the blank separator) stays untouched.
# -------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$p3 = $d.Paragraphs(3)

$listRange = $d.Range($p1.Range.Start, $p3.Range.End)
$listRange.set_Style("List Paragraph")
$listRange.ListFormat.ApplyBulletDefault()

$p2.Range.ListFormat.RemoveNumbers()
$p2.Range.set_Style("Normal")

# -------------------------------------------------------------------
# Make the generated "List Paragraph" style match Word's built-in
# definition (indentation + contextual spacing + ui priority).
# -------------------------------------------------------------------
$listParagraphStyle = $d.Styles("List Paragraph")
$listParagraphStyle.Priority = 34
$listParagraphStyle.ParagraphFormat.LeftIndent = 36
$listParagraphStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

Write-Output "edit complete"
